# Refresh the cryptocurrency price/volume snapshot (Price column D,
# Volume(1h) column E) for rows 2-51 on Sheet1, per the latest pull from
# the GitHub Actions scheduled job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "67.627.71"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -0.79%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.483.53"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -1.41%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.04%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "592.49"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -1.88%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "179.20"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -1.98%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  +2.03%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  -0.01%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "3.483.43"; ForceText = $false },
    @{ Cell = "E9"; Value = "  -1.38%  "; ForceText = $false },
    @{ Cell = "E10"; Value = "  -1.63%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "6.99"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -2.78%  "; ForceText = $false },
    @{ Cell = "E12"; Value = "  -3.33%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "4.091.43"; ForceText = $false },
    @{ Cell = "E13"; Value = "  -1.25%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "32.54"; ForceText = $true },
    @{ Cell = "E14"; Value = "  +0.45%  "; ForceText = $false },
    @{ Cell = "E15"; Value = "  -2.57%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "67.614.58"; ForceText = $false },
    @{ Cell = "E16"; Value = "  -0.68%  "; ForceText = $false },
    @{ Cell = "E17"; Value = "  -2.64%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "3.481.83"; ForceText = $false },
    @{ Cell = "E18"; Value = "  -1.40%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "6.15"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -4.03%  "; ForceText = $false },
    @{ Cell = "E20"; Value = "  -3.44%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "389.65"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -2.97%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "7.94"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -1.70%  "; ForceText = $false },
    @{ Cell = "E23"; Value = "  +1.45%  "; ForceText = $false },
    @{ Cell = "E24"; Value = "  +0.35%  "; ForceText = $false },
    @{ Cell = "E25"; Value = "  -2.26%  "; ForceText = $false },
    @{ Cell = "E26"; Value = "  -2.00%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "0.0000123"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -1.17%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "10.15"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -4.03%  "; ForceText = $false },
    @{ Cell = "E29"; Value = "  -1.17%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E30"; Value = "  +0.37%  "; ForceText = $false },
    @{ Cell = "E31"; Value = "  -5.07%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "24.76"; ForceText = $true },
    @{ Cell = "E32"; Value = "  +2.92%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "2.04"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -2.07%  "; ForceText = $false },
    @{ Cell = "E34"; Value = "  -5.16%  "; ForceText = $false },
    @{ Cell = "E35"; Value = "  -3.50%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  -0.09%  "; ForceText = $false },
    @{ Cell = "E37"; Value = "  -4.81%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "161.07"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -1.23%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "0.890"; ForceText = $true },
    @{ Cell = "E39"; Value = "  +0.92%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "28.09"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +5.44%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "1.85"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -4.74%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  -4.75%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "6.62"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -6.57%  "; ForceText = $false },
    @{ Cell = "E44"; Value = "  -4.49%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "0.0713"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -3.78%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "2.727.84"; ForceText = $false },
    @{ Cell = "E46"; Value = "  -6.71%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "26.00"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -3.49%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "41.52"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -2.36%  "; ForceText = $false },
    @{ Cell = "E49"; Value = "  -2.86%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "332.40"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -5.78%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "1.04"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -3.97%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Some "Price" values are plain decimals (e.g. "592.49") that Excel's
        # type-inference would otherwise silently convert to a Number,
        # dropping formatting / precision. Force them to stay text, matching
        # the rest of the column (which holds thousand-dotted strings like
        # "67.627.71" that can never parse as numbers anyway), then clear
        # the now-unneeded explicit "Text" number format so the cell's style
        # matches its untouched neighbours again.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
